{"js": "// Add the GitHub link as a new run at the end of the paragraph that\n// introduces it (\"\u041f\u043e\u0441\u0438\u043b\u0430\u043d\u043d\u044f \u043d\u0430 git hub: \").\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst needle = \"git hub\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text && para.text.indexOf(needle) !== -1) {\n    target = para;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find the \"\u041f\u043e\u0441\u0438\u043b\u0430\u043d\u043d\u044f \u043d\u0430 git hub:\" paragraph');\n}\n\ntarget.insertText(\n  \"https://github.com/Oleksiy2003/Labs/blob/main/Lab4/Lab4.3/Lab4.3/main.cpp\",\n  Word.InsertLocation.end\n);\nawait context.sync();\n", "ps1": "# Add the GitHub link as a new run at the end of the paragraph that\n# introduces it (\"\u041f\u043e\u0441\u0438\u043b\u0430\u043d\u043d\u044f \u043d\u0430 git hub: \").\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$found = $range.Find.Execute(\"git hub:\")\n\nif ($found) {\n    $range.Expand(4) | Out-Null        # wdParagraph - grow to the whole paragraph\n    $range.MoveEnd(1, -1) | Out-Null   # wdCharacter - back off the paragraph mark\n    $range.InsertAfter(\"https://github.com/Oleksiy2003/Labs/blob/main/Lab4/Lab4.3/Lab4.3/main.cpp\") | Out-Null\n}\n"}
